$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove now-unused rows 42:53 (table shrinks from 53 to 41 rows)
$ws.Rows("42:53").Delete()

# Update header-adjacent data rows (2-41) with final resolved values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C2").Value = '1. Object Naming'
$ws.Range("D2").Value = 'Title'
$ws.Range("E2").Value = 'Capitalizar as inicais de nomes próprios e da primeira palavra, para outros termos use letras minúsculas'
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 1215
$ws.Range("H2").Value = 1215

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C3").Value = '1. Object Naming'
$ws.Range("D3").Value = 'Title'
$ws.Range("E3").Value = 'Evitar abreviações'
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 1215
$ws.Range("H3").Value = 1215

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C4").Value = '1. Object Naming'
$ws.Range("D4").Value = 'Title'
$ws.Range("E4").Value = 'Não pode ficar vazio'
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1215
$ws.Range("H4").Value = 1215

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C5").Value = '1. Object Naming'
$ws.Range("D5").Value = 'Title'
$ws.Range("E5").Value = 'Não utilizar artigos'
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1215
$ws.Range("H5").Value = 1215

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C6").Value = '1. Object Naming'
$ws.Range("D6").Value = 'Work Type'
$ws.Range("E6").Value = 'Capitalizar as inicais de nomes próprios e da primeira palavra, para outros termos use letras minúsculas'
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1214
$ws.Range("H6").Value = 1215

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C7").Value = '1. Object Naming'
$ws.Range("D7").Value = 'Work Type'
$ws.Range("E7").Value = 'Evitar abreviações'
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 1215
$ws.Range("H7").Value = 1215

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C8").Value = '1. Object Naming'
$ws.Range("D8").Value = 'Work Type'
$ws.Range("E8").Value = 'Fazer uso de vocabulário controlado'
$ws.Range("F8").Value = 1215
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 1215

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C9").Value = '1. Object Naming'
$ws.Range("D9").Value = 'Work Type'
$ws.Range("E9").Value = 'Não pode ficar vazio'
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 1215
$ws.Range("H9").Value = 1215

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C10").Value = '1. Object Naming'
$ws.Range("D10").Value = 'Work Type'
$ws.Range("E10").Value = 'Não usar pontuação, exceto hífen'
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1214
$ws.Range("H10").Value = 1215

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C11").Value = '1. Object Naming'
$ws.Range("D11").Value = 'Work Type'
$ws.Range("E11").Value = 'Usar o mesmo idioma do catálogo'
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 1215
$ws.Range("H11").Value = 1215

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C12").Value = '1. Object Naming'
$ws.Range("D12").Value = 'Work Type'
$ws.Range("E12").Value = 'Usar singular'
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 1215
$ws.Range("H12").Value = 1215

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C13").Value = '2. Creator Information'
$ws.Range("D13").Value = 'Creator'
$ws.Range("E13").Value = 'Evitar abreviações'
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 1215
$ws.Range("H13").Value = 1215

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C14").Value = '2. Creator Information'
$ws.Range("D14").Value = 'Creator'
$ws.Range("E14").Value = 'Fazer uso de vocabulário controlado'
$ws.Range("F14").Value = 1215
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 1215

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C15").Value = '2. Creator Information'
$ws.Range("D15").Value = 'Creator'
$ws.Range("E15").Value = 'Não pode ficar vazio'
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 1215
$ws.Range("H15").Value = 1215

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C16").Value = '3. Physical Characteristics'
$ws.Range("D16").Value = 'Materials and Techniques'
$ws.Range("E16").Value = 'Capitalizar as inicais de nomes próprios e da primeira palavra, para outros termos use letras minúsculas'
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 1215
$ws.Range("H16").Value = 1215

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C17").Value = '3. Physical Characteristics'
$ws.Range("D17").Value = 'Materials and Techniques'
$ws.Range("E17").Value = 'Evitar abreviações'
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 1215
$ws.Range("H17").Value = 1215

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C18").Value = '3. Physical Characteristics'
$ws.Range("D18").Value = 'Materials and Techniques'
$ws.Range("E18").Value = 'Fazer uso de vocabulário controlado'
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 1215
$ws.Range("H18").Value = 1215

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C19").Value = '3. Physical Characteristics'
$ws.Range("D19").Value = 'Materials and Techniques'
$ws.Range("E19").Value = 'Não pode ficar vazio'
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 1215
$ws.Range("H19").Value = 1215

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C20").Value = '3. Physical Characteristics'
$ws.Range("D20").Value = 'Materials and Techniques'
$ws.Range("E20").Value = 'Usar o mesmo idioma do catálogo'
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 1215
$ws.Range("H20").Value = 1215

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C21").Value = '3. Physical Characteristics'
$ws.Range("D21").Value = 'Materials and Techniques'
$ws.Range("E21").Value = 'Usar singular'
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 1215
$ws.Range("H21").Value = 1215

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C22").Value = '3. Physical Characteristics'
$ws.Range("D22").Value = 'Physical Description'
$ws.Range("E22").Value = 'Fazer uso de vocabulário controlado'
$ws.Range("F22").Value = 1215
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 1215

$ws.Range("A23").Value = 21
$ws.Range("B23").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C23").Value = '4. Stylistic, Cultural, and Chronological Information'
$ws.Range("D23").Value = 'Date'
$ws.Range("E23").Value = 'Anos com menos que 4 digitos, inserir 0 a esquerda'
$ws.Range("F23").Value = 72
$ws.Range("G23").Value = 1143
$ws.Range("H23").Value = 1215

$ws.Range("A24").Value = 22
$ws.Range("B24").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C24").Value = '4. Stylistic, Cultural, and Chronological Information'
$ws.Range("D24").Value = 'Date'
$ws.Range("E24").Value = 'Capitalizar as inicais de nomes próprios e da primeira palavra, para outros termos use letras minúsculas'
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 1213
$ws.Range("H24").Value = 1215

$ws.Range("A25").Value = 23
$ws.Range("B25").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C25").Value = '4. Stylistic, Cultural, and Chronological Information'
$ws.Range("D25").Value = 'Date'
$ws.Range("E25").Value = 'Não pode ficar vazio'
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 1215
$ws.Range("H25").Value = 1215

$ws.Range("A26").Value = 24
$ws.Range("B26").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C26").Value = '4. Stylistic, Cultural, and Chronological Information'
$ws.Range("D26").Value = 'Date'
$ws.Range("E26").Value = 'Não utilizar apostrofo'
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 1215
$ws.Range("H26").Value = 1215

$ws.Range("A27").Value = 25
$ws.Range("B27").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C27").Value = '4. Stylistic, Cultural, and Chronological Information'
$ws.Range("D27").Value = 'Date'
$ws.Range("E27").Value = 'Seguir padrão para registro de hora, minutos e segundos'
$ws.Range("F27").Value = 1215
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 1215

$ws.Range("A28").Value = 26
$ws.Range("B28").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C28").Value = '4. Stylistic, Cultural, and Chronological Information'
$ws.Range("D28").Value = 'Date'
$ws.Range("E28").Value = 'Seguir padrão pra registro de dia, mês e ano de data'
$ws.Range("F28").Value = 1215
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 1215

$ws.Range("A29").Value = 27
$ws.Range("B29").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C29").Value = '4. Stylistic, Cultural, and Chronological Information'
$ws.Range("D29").Value = 'Date'
$ws.Range("E29").Value = 'Usar o mesmo idioma do catálogo'
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 1215
$ws.Range("H29").Value = 1215

$ws.Range("A30").Value = 28
$ws.Range("B30").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C30").Value = '4. Stylistic, Cultural, and Chronological Information'
$ws.Range("D30").Value = 'Date'
$ws.Range("E30").Value = 'Use traço para separar intervalo de anos'
$ws.Range("F30").Value = 1156
$ws.Range("G30").Value = 59
$ws.Range("H30").Value = 1215

$ws.Range("A31").Value = 29
$ws.Range("B31").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C31").Value = '5. Location and Geography'
$ws.Range("D31").Value = 'Creation Location'
$ws.Range("E31").Value = 'Capitalizar as inicais de nomes próprios e da primeira palavra, para outros termos use letras minúsculas'
$ws.Range("F31").Value = 118
$ws.Range("G31").Value = 1097
$ws.Range("H31").Value = 1215

$ws.Range("A32").Value = 30
$ws.Range("B32").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C32").Value = '5. Location and Geography'
$ws.Range("D32").Value = 'Creation Location'
$ws.Range("E32").Value = 'Evitar abreviações'
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 1215
$ws.Range("H32").Value = 1215

$ws.Range("A33").Value = 31
$ws.Range("B33").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C33").Value = '5. Location and Geography'
$ws.Range("D33").Value = 'Creation Location'
$ws.Range("E33").Value = 'Fazer uso de vocabulário controlado'
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 1215
$ws.Range("H33").Value = 1215

$ws.Range("A34").Value = 32
$ws.Range("B34").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C34").Value = '5. Location and Geography'
$ws.Range("D34").Value = 'Creation Location'
$ws.Range("E34").Value = 'Usar o mesmo idioma do catálogo'
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 1215
$ws.Range("H34").Value = 1215

$ws.Range("A35").Value = 33
$ws.Range("B35").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C35").Value = '7. Class'
$ws.Range("D35").Value = 'Class'
$ws.Range("E35").Value = 'Evitar abreviações'
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 1215
$ws.Range("H35").Value = 1215

$ws.Range("A36").Value = 34
$ws.Range("B36").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C36").Value = '7. Class'
$ws.Range("D36").Value = 'Class'
$ws.Range("E36").Value = 'Fazer uso de vocabulário controlado'
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 1215
$ws.Range("H36").Value = 1215

$ws.Range("A37").Value = 35
$ws.Range("B37").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C37").Value = '7. Class'
$ws.Range("D37").Value = 'Class'
$ws.Range("E37").Value = 'Usar singular'
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 1215
$ws.Range("H37").Value = 1215

$ws.Range("A38").Value = 36
$ws.Range("B38").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C38").Value = '8. Description'
$ws.Range("D38").Value = 'Description'
$ws.Range("E38").Value = 'Capitalizar as inicais de nomes próprios e da primeira palavra, para outros termos use letras minúsculas'
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = 1214
$ws.Range("H38").Value = 1215

$ws.Range("A39").Value = 37
$ws.Range("B39").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C39").Value = '8. Description'
$ws.Range("D39").Value = 'Description'
$ws.Range("E39").Value = 'Evitar abreviações'
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 1215
$ws.Range("H39").Value = 1215

$ws.Range("A40").Value = 38
$ws.Range("B40").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C40").Value = '8. Description'
$ws.Range("D40").Value = 'Description'
$ws.Range("E40").Value = 'Fazer uso de vocabulário controlado'
$ws.Range("F40").Value = 1215
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 1215

$ws.Range("A41").Value = 39
$ws.Range("B41").Value = 'mhn_-_moedas-de-ouro.csv'
$ws.Range("C41").Value = '8. Description'
$ws.Range("D41").Value = 'Description'
$ws.Range("E41").Value = 'Usar o mesmo idioma do catálogo'
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 1215
$ws.Range("H41").Value = 1215
